# Insert a new data row at row 23 (shifting existing rows 23:87 down to 24:88)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("23:23").Insert()

$ws.Cells.Item(23, 1).Value2  = 1
$ws.Cells.Item(23, 2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(23, 3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(23, 4).Value2  = 44742
$ws.Cells.Item(23, 5).Value2  = 15
$ws.Cells.Item(23, 6).Value2  = "Fruta"
$ws.Cells.Item(23, 7).Value2  = 100102
$ws.Cells.Item(23, 8).Value2  = "Cítricos"
$ws.Cells.Item(23, 9).Value2  = 100102005
$ws.Cells.Item(23, 10).Value2 = "Naranja"
$ws.Cells.Item(23, 11).Value2 = "Navel"
$ws.Cells.Item(23, 12).Value2 = "Tercera"
$ws.Cells.Item(23, 13).Value2 = 300
$ws.Cells.Item(23, 14).Value2 = 600
$ws.Cells.Item(23, 15).Value2 = 650
$ws.Cells.Item(23, 16).Value2 = 625
$ws.Cells.Item(23, 17).Value2 = "`$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(23, 18).Value2 = "Región de Coquimbo"
$ws.Cells.Item(23, 19).Value2 = 625
$ws.Cells.Item(23, 20).Value2 = 1
